# Bitacora.xlsx update: add a new log entry for "Creacion de controller
# expenses y varios metodos" to row 27 (mirrors the style of row 17, the
# most recent entry that used the same "Controllers" topic), and move the
# active selection to F30 as left by the author after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 / 28 are a pre-formatted but still-empty pair of merged placeholder
# rows (C27:C28, D27:D28, E27:E28). Copy the cell formatting used by the
# previous "Controllers" entry (row 17) into row 27 so the new entry gets
# the same borders/fonts/number-format, then stamp in the new values.
$ws.Range("C17:E17").Copy()
$ws.Range("C27:E27").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 17 also has a custom row height (15pt) that differs from the default
# placeholder row height - match it.
$ws.Rows.Item(27).RowHeight = $ws.Rows.Item(17).RowHeight

# Fill in the new bitacora entry.
$ws.Range("C27").Value2 = "Controllers"
$ws.Range("D27").Value2 = 44741
$ws.Range("E27").Value2 = "Creacion de controller expenses y varios metodos"

# Leave the selection where the author left it after typing the new row.
$ws.Range("F30").Select()
